# 04/12 - Navigation On Company and TestData File updates
#
# Renames the "TBD" placeholder sheets to "OnGoing" and replaces their
# stale sample data with the current environment / partner links.

$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("Environments_TBD")
$wsEnv.Name = "Environments_OnGoing"

$wsUsers = $wb.Worksheets.Item("Users_TBD")
$wsUsers.Name = "Users_OnGoing"

# --- Users_OnGoing: drop the stale rows, keep a handful of blank,
#     hyperlink-styled placeholder cells -----------------------------------
$wsUsers.Hyperlinks.Delete()
$wsUsers.Range("A2:B13").Clear()

$wsUsers.Range("B2").Style = "Hyperlink"
$wsUsers.Range("B3").Style = "Hyperlink"
$wsUsers.Range("B6").Style = "Hyperlink"
$wsUsers.Range("B7").Style = "Hyperlink"
$wsUsers.Range("B9").Style = "Hyperlink"

$wsUsers.Range("B26").Select()

# --- Environments_OnGoing: replace old env rows with the two live links --
$wsEnv.Hyperlinks.Delete()
$wsEnv.Range("A2:B21").Clear()

$wsEnv.Range("A2").Value = "DirectSales"
$wsEnv.Range("B2").Value = "https://proximus--prxitt.my.salesforce.com"

$wsEnv.Range("A3").Value = "PartnersCommunity"
$wsEnv.Range("B3").Value = "https://prxitt-proximus.cs127.force.com/SalesforceforPartners"

$wsEnv.Hyperlinks.Add($wsEnv.Range("B2"), "https://proximus--prxitt.my.salesforce.com") | Out-Null
$wsEnv.Hyperlinks.Add($wsEnv.Range("B3"), "https://prxitt-proximus.cs127.force.com/SalesforceforPartners") | Out-Null

$wsEnv.Range("B2").Style = "Hyperlink"
$wsEnv.Range("B3").Style = "Hyperlink"
$wsEnv.Range("B4").Style = "Hyperlink"
$wsEnv.Range("B5").Style = "Hyperlink"
$wsEnv.Range("B6").Style = "Hyperlink"

$wsEnv.Range("A4").Select()
